$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.961.23'
$ws.Range("E2").Value = '  +5.70%  '
$ws.Range("D3").Value = '3.646.45'
$ws.Range("E3").Value = '  +16.50%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.61%  '
$ws.Range("D7").Value = '3.645.32'
$ws.Range("E7").Value = '  +16.62%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +3.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.497'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000253'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.41%  '
$ws.Range("D15").Value = '4.257.98'
$ws.Range("E15").Value = '  +16.55%  '
$ws.Range("D16").Value = '70.944.67'
$ws.Range("E16").Value = '  +5.81%  '
$ws.Range("D17").Value = '3.653.78'
$ws.Range("E17").Value = '  +16.72%  '
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '513.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +16.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.741'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.37%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("E31").Value = '  +17.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.44'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.116'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.65%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.10'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.88%  '
$ws.Range("E37").Value = '  +6.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.344'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.01%  '
$ws.Range("E41").Value = '  +4.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.93%  '
$ws.Range("D43").Value = '3.126.12'
$ws.Range("E43").Value = '  +11.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '417.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.01%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +14.68%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0367'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.11%  '
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.61%  '
